# Auto update: 2025-05-20 20:38:52
# Reshuffle Company Name / Company Number / Category values across rows
# (columns C-G stay anchored to their row; this mirrors an upstream re-sort
# of the underlying source data for Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    ,@(2, 'DAVIDSON CAPITAL HOLDINGS LTD', 'SC849117', 'Capital')
    ,@(3, 'ST GEORGE CAPITAL (LAND) LIMITED', '16462880', 'Capital')
    ,@(4, 'AFROSCOT VENTURES LTD', '16462878', 'Ventures')
    ,@(6, 'SAMVIV PARTNERS LTD', '16460672', 'Partners')
    ,@(7, 'T GILPIN PHYSIO CONSULTANCY LTD', '16460503', 'LP')
    ,@(8, '4D CAPITAL PROPCO (44) LIMITED', '16461269', 'Capital')
    ,@(9, 'THE REEL MED LLP', 'OC456780', 'LP')
    ,@(10, 'PONGPONG MALATANG LTD', '16458077', 'GP')
    ,@(11, 'KNOTT INVESTMENTS LIMITED', '16458684', 'Investments')
    ,@(12, 'KC INVESTMENTS & TRADING LIMITED', '16456642', 'Investments')
    ,@(13, 'JJOHN INVESTMENTS LIMITED', '16456276', 'Investments')
    ,@(14, 'ECHO VENTURES GROUP LIMITED', '16455744', 'Ventures')
    ,@(16, 'MUSICROOTS LTD', '16455514', 'SIC')
    ,@(17, 'ARISSA INVESTMENTS LIMITED', '16455197', 'Investments')
    ,@(18, 'TALKSGPT AI LTD', '16455313', 'GP')
    ,@(21, 'TUERNER IMMIGRATION LLP', 'OC456770', 'LP')
    ,@(23, 'GROWTHFORGE MANAGEMENT LLP', 'OC456769', 'LP')
    ,@(24, 'GOLDEN VENTURES LONDON LTD', '16452104', 'Ventures')
    ,@(25, 'ALDABBOUS UK INVESTMENTS LTD', '16453476', 'Investments')
    ,@(26, 'CAMBRIDGE SOCIAL INVESTMENTS LIMITED', '16453466', 'Investments')
    ,@(27, 'CAPITAL & CENTRIC (SYNCHRONICITY) LTD', '16453716', 'Capital')
    ,@(28, 'GULF TRADE AND INVESTMENT ADVANTAGES JOINT PARTNERSHIP LTD', '16453733', 'Partners')
    ,@(29, 'FROST CAPITAL LTD', '16450073', 'Capital')
    ,@(31, 'ASSET CAPITAL 44 OPCO LIMITED', '16449512', 'Capital')
    ,@(32, 'ATHENA PARTNERSHIP LTD', '16449517', 'Partners')
)

foreach ($u in $updates) {
    $row = $u[0]
    $companyName = $u[1]
    $companyNumber = $u[2]
    $category = $u[3]

    $ws.Cells.Item($row, 1).Value = $companyName

    # Company numbers that are purely digits (e.g. "16462880") would
    # otherwise be auto-coerced to a numeric cell by Excel; prefix with
    # a quote so they stay text, exactly like the other untouched rows
    # (e.g. "SC849118", "OC456771") that are already stored as text.
    if ($companyNumber -match "^[0-9]+$") {
        $ws.Cells.Item($row, 2).Value = "'" + $companyNumber
    } else {
        $ws.Cells.Item($row, 2).Value = $companyNumber
    }

    $ws.Cells.Item($row, 8).Value = $category
}
